$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 50
$ws1.Range("F3").Value = 253
$ws1.Range("F7").Value = 614
$ws1.Range("F8").Value = 226
$ws1.Range("F11").Value = 158
$ws1.Range("F12").Value = 714
$ws1.Range("F13").Value = 96
$ws1.Range("F14").Value = 1837
$ws1.Range("F15").Value = 371
$ws1.Range("F16").Value = 3924
$ws1.Range("F18").Value = 498
$ws1.Range("F20").Value = 61
$ws1.Range("F21").Value = 149

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 23
$ws2.Range("F7").Value = 480
$ws2.Range("F17").Value = 31
$ws2.Range("F21").Value = 24

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 331

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

# Step 1: apply F-column updates that exist BEFORE the duplicate row is removed
$ws4.Range("F2").Value = 50
$ws4.Range("F4").Value = 331
$ws4.Range("F7").Value = 253
$ws4.Range("F10").Value = 23
$ws4.Range("F12").Value = 480

# Step 2: remove the duplicate "神山羊2024巡演ENCOUNTER" row (row 13); rows below shift up by one
$ws4.Rows.Item(13).Delete()

# Step 3: fix the literal running-index values in column A for every row that shifted up
for ($r = 13; $r -le 45; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

# Step 4: apply F-column updates addressed at their NEW (post-shift) row numbers
$ws4.Range("F18").Value = 614
$ws4.Range("F19").Value = 226
$ws4.Range("F23").Value = 158
$ws4.Range("F26").Value = 714
$ws4.Range("F27").Value = 96
$ws4.Range("F29").Value = 1837
$ws4.Range("F31").Value = 3924
$ws4.Range("F34").Value = 498
$ws4.Range("F36").Value = 61
$ws4.Range("F38").Value = 149
$ws4.Range("F41").Value = 31
$ws4.Range("F45").Value = 24

